$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.022.35"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.643.68"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.84"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0639"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "1.663.46"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.43"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "26.044.46"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "194.12"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.14"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.24"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.906"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "1.130.26"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.539"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.45"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.04"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.797"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.74"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0950"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  +3.08%  "
